$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values for the affected rows.
# These new values reflect the "repull data, push all data, mean calculation" update.
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -5
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = 7
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = -4
$ws.Range("F13").Value = -3
$ws.Range("F14").Value = -1
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = -4
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = -2
$ws.Range("F20").Value = -3
